$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 00:55:58"
$wsZhCn.Range("H2").Value = "2016-03-21 00:56:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 00:56:02"
$wsDeDe.Range("H2").Value = "2016-03-21 00:56:23"
